$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.011.14"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -0.28%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.931.11"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.95%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'354.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.63%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'107.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -5.45%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.564"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +1.43%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +0.17%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.620"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'38.20"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -3.67%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +1.10%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0858"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -0.57%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'19.15"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -3.06%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'3.404.98"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +1.37%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'7.68"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -0.53%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'2.908.75"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.00%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'0.971"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -1.19%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'52.019.50"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -0.35%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'3.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +3.13%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'7.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -1.20%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'13.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -2.96%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.0₃0974"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -0.17%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'69.90"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -1.87%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'266.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -1.22%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  -1.16%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.175"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -3.59%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'26.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +0.42%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'7.59"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +12.39%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  +1.84%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'10.33"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -2.76%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'36.65"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -2.12%  "
$ws.Range("E32").ClearFormats()
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'2.18"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -4.20%  "
$ws.Range("E33").ClearFormats()
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'5.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -3.59%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'52.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -2.05%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.0434"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -3.63%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +0.05%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'3.17"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -4.20%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  -2.42%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'17.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -4.70%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'2.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -1.18%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.118"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +0.88%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'23.02"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.09%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'118.79"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.02%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'  -0.83%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'2.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -3.12%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'2.118.65"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -2.55%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'3.37"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -3.71%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'3.237.18"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +1.47%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.241"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -7.96%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'  -0.17%  "
$ws.Range("E51").ClearFormats()
